$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 292, shifting rows 292:403 down to 293:404.
$ws.Rows(292).Insert()

# Populate the newly inserted row 292 with its data.
$ws.Cells.Item(292, 1).Value = 9
$ws.Cells.Item(292, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(292, 3).Value = "Metropolitana"
$ws.Cells.Item(292, 4).Value = 45229
$ws.Cells.Item(292, 5).Value = 13
$ws.Cells.Item(292, 6).Value = 100112001
$ws.Cells.Item(292, 7).Value = "Berenjena"
$ws.Cells.Item(292, 8).Value = "Sin especificar"
$ws.Cells.Item(292, 9).Value = "Primera"
$ws.Cells.Item(292, 10).Value = 160
$ws.Cells.Item(292, 11).Value = 10000
$ws.Cells.Item(292, 12).Value = 11000
$ws.Cells.Item(292, 13).Value = 10500
$ws.Cells.Item(292, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(292, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(292, 16).Value = 210
$ws.Cells.Item(292, 17).Value = 50
$ws.Cells.Item(292, 18).Value = "Hortaliza"
